# Split the opening sentence into two runs wrapped with grammar proofErr
# markers, then add a new "TEST TEST TEST" paragraph (middle word wrapped
# with spell-check proofErr markers) right before the existing
# bookmarkStart/bookmarkEnd ("_GoBack") pair.

$d = $word.ActiveDocument

# The whole body story ("Hey I'm on the doc right now" + the _GoBack
# bookmark) lives in a single paragraph before this edit; InsertXML on
# the full Content range lets us replace it with the exact target markup
# (including the w:proofErr grammar/spelling markers that Word's
# proofing pass would have inserted) in one shot, re-declaring the
# _GoBack bookmark on the new second paragraph so it isn't lost.
$rng = $d.Content

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
          '<w:p>' + `
            '<w:proofErr w:type="gramStart"/>' + `
            '<w:r><w:t>Hey</w:t></w:r>' + `
            '<w:proofErr w:type="gramEnd"/>' + `
            '<w:r><w:t xml:space="preserve"> I’m on the doc right now</w:t></w:r>' + `
          '</w:p>' + `
          '<w:p>' + `
            '<w:r><w:t xml:space="preserve">TEST </w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:t>TEST</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:t xml:space="preserve"> TEST</w:t></w:r>' + `
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
            '<w:bookmarkEnd w:id="0"/>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'

$rng.InsertXML($xml)
